$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update column C (row 2 through row 13) from 45221 to 45224 ("Förändrad" date column)
$ws.Range("C2:C13").Value = 45224
